# Applies the "Donttest" column addition to the checklist worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column F: header + per-function values.
$ws.Range("F4").Value = "Donttest"
$ws.Range("F5").Value = "YES"
$ws.Range("F6").Value = "NULL"
$ws.Range("F7").Value = "YES"
$ws.Range("F8").Value = "YES"
$ws.Range("F9").Value = "YES"
$ws.Range("F10").Value = "YES"
$ws.Range("F11").Value = "YES"
$ws.Range("F12").Value = "YES"
$ws.Range("F13").Value = "YES"
$ws.Range("F14").Value = "YES"

# Update the active selection to match the post-edit workbook state.
$ws.Range("F15").Select()
